$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'What is the process of data analysis?'
$ws.Range("B2").Value = 'Data analysis is a systematic process that involves inspecting, cleaning, transforming, and modeling data to extract useful information and support decision-making. The process begins with data collection, followed by data cleaning to eliminate errors or inconsistencies. Data transformation is done to convert data into a suitable format for analysis. Then, various analytical techniques, such as statistical analysis or data mining, are applied to explore patterns, relationships, or trends in the data. Finally, the results are interpreted and presented in a meaningful way to inform decision-making.'
$ws.Range("C2").Value = 'Technical'

$ws.Range("A3").Value = 'Can you describe the concept of data visualization?'
$ws.Range("B3").Value = 'Data visualization is the presentation of data in a graphical or visual format that facilitates understanding and interpretation. It involves the use of charts, graphs, maps, or other visual elements to represent data patterns, relationships, or trends. The primary purpose of data visualization is to make complex data more accessible, allowing users to identify patterns or outliers quickly and make informed decisions. By presenting data visually, it becomes easier to communicate insights, compare values, and understand the significance of the data. Effective data visualization enhances data analysis by providing a clear and concise representation of information.'
$ws.Range("C3").Value = 'Technical'

$ws.Range("A4").Value = 'Describe a time when you had to work under tight deadlines for a data analysis project.'
$ws.Range("B4").Value = 'When answering this question, it is important to highlight your ability to work efficiently and effectively under pressure. Discuss a specific situation where you had to complete a data analysis project within a tight timeframe. Start by explaining the deadline and the importance of the project. Then, outline the steps you took to manage your time and prioritize tasks. Mention any strategies or techniques you used to streamline the data analysis process without compromising accuracy. Finally, share the outcome of the project and how you successfully met the deadline.'
$ws.Range("C4").Value = 'Behavioral'

$ws.Range("A5").Value = 'Explain the concept of data normalization.'
$ws.Range("B5").Value = 'Data normalization is the process of organizing and structuring data to reduce redundancy and anomalies. It involves eliminating data duplication and ensuring data consistency by adhering to specific rules or normalization forms. The purpose of data normalization is to minimize data integrity issues, enhance data storage efficiency, and streamline data analysis. By normalizing data, you can eliminate data inconsistencies and improve the accuracy and reliability of analysis results. Common normalization techniques include splitting data into tables, creating relationships between tables, and applying normalization rules such as the first normal form (1NF), second normal form (2NF), and third normal form (3NF).'
$ws.Range("C5").Value = 'Technical'

$ws.Range("A6").Value = 'Describe a situation where you had to present complex data findings to a non-technical audience.'
$ws.Range("B6").Value = 'When answering this question, focus on your ability to communicate complex information effectively. Start by explaining the context and the audience''s level of technical knowledge. Outline the key insights or findings from your data analysis and consider how to simplify and present the information in a clear and concise manner. Use visual aids, such as charts or graphs, to support your explanation and make the data more accessible. Additionally, emphasize your communication skills, including the use of plain language, storytelling techniques, or analogies to help the audience understand the data. Finally, mention any feedback or questions received from the audience and how you addressed them.'
$ws.Range("C6").Value = 'Behavioral'

$ws.Range("A7").Value = 'What is the role of exploratory data analysis in the data analysis process?'
$ws.Range("B7").Value = 'Exploratory data analysis (EDA) is a critical step in the data analysis process. Its primary role is to gain insights and understand the data before applying more complex analysis techniques. During EDA, you explore the data to identify patterns, trends, or outliers. This involves the use of summary statistics, visualizations, and basic analytical techniques. EDA helps in understanding data quality, detecting data errors or missing values, and uncovering relationships or associations between variables. By gaining a deeper understanding of the data through EDA, you can make informed decisions about which analysis methods are appropriate and how to interpret the results accurately.'
$ws.Range("C7").Value = 'Technical'

$ws.Rows("8:11").Delete()
